$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B30").Value = 40
$ws.Range("D30").Value = 19
$ws.Range("E30").Value = 'This dataset focuses on population change by regions. The data is collected at the yearly level and covers the territory types of State and Region. The dataset includes the following columns: "Roky", ...'
$ws.Range("F30").Value = 'WREG01AT2'
$ws.Range("G30").Value = 0.3456027507781982
$ws.Range("H30").Value = 0.5946824685399295
$ws.Range("I30").Value = 0.002560871
$ws.Range("B31").Value = 41
$ws.Range("D31").Value = 13
$ws.Range("E31").Value = '"This dataset focuses on the main data on population and vital statistics, specifically relative figures. The data is collected at the yearly level and covers two types of territories: the state and r...'
$ws.Range("F31").Value = 'OBY01BT02'
$ws.Range("G31").Value = 0.3427191972732544
$ws.Range("H31").Value = 0.5908731159969555
$ws.Range("I31").Value = 0.0052404497
$ws.Range("B32").Value = 47
$ws.Range("D32").Value = 12
$ws.Range("E32").Value = '"This dataset focuses on population and population change since 1785, presented in absolute figures. The data is collected at the yearly level and covers the territory type "State". The dataset includ...'
$ws.Range("F32").Value = 'OBY01CT01'
$ws.Range("G32").Value = 0.3326067924499512
$ws.Range("H32").Value = 0.5759289386249279
$ws.Range("I32").Value = 0.0067449897
$ws.Range("B33").Value = 49
$ws.Range("D33").Value = 47
$ws.Range("E33").Value = 'This dataset focuses on occupied houses with dwellings categorized by the number of dwellings and regions. The data is collected at the "Census Year" level and covers the territory types "State" and "...'
$ws.Range("F33").Value = 'SLD018T02'
$ws.Range("G33").Value = 0.3284686803817749
$ws.Range("H33").Value = 0.5758202018419979
$ws.Range("I33").Value = 0.000013315579
$ws.Range("B34").Value = 42
$ws.Range("D34").Value = 35
$ws.Range("E34").Value = '"This dataset focuses on the population distribution by sex, marital status, and region. The data is collected at the ''Census Year'' level and covers the territory types of ''State'' and ''Region''. The da...'
$ws.Range("F34").Value = 'SLD001T02'
$ws.Range("G34").Value = 0.3406649827957153
$ws.Range("H34").Value = 0.5751797020827104
$ws.Range("I34").Value = 0.00021318757
$ws.Range("D35").Value = 20
$ws.Range("E35").Value = '"This dataset focuses on the main data regarding population and vital statistics by municipalities. The data is collected at the yearly level and covers regions and municipalities. The dataset include...'
$ws.Range("F35").Value = 'OBY01B01T01'
$ws.Range("G35").Value = 0.3311309218406677
$ws.Range("H35").Value = 0.5717562759906273
$ws.Range("I35").Value = 0.0025409926
$ws.Range("B36").Value = 43
$ws.Range("D36").Value = 34
$ws.Range("E36").Value = '"This dataset focuses on fertility rates of women by five-year age groups. The data is collected at the yearly level and covers the territory types of State and Region. The dataset includes the follow...'
$ws.Range("F36").Value = 'OBY03DT02'
$ws.Range("G36").Value = 0.3395309448242188
$ws.Range("H36").Value = 0.5713558661567641
$ws.Range("I36").Value = 0.00022693386
$ws.Range("B37").ClearContents()
$ws.Range("D37").Value = 43
$ws.Range("E37").Value = 'This dataset focuses on "Births - selected summary data". The data is collected at the "Year" level and covers the territory types "State" and "Region". The dataset includes the following columns: "Ro...'
$ws.Range("F37").Value = 'OBY03T01'
$ws.Range("G37").ClearContents()
$ws.Range("H37").Value = 0.5703342306040304
$ws.Range("I37").Value = 0.00003676463
$ws.Range("B38").Value = 56
$ws.Range("D38").Value = 39
$ws.Range("E38").Value = 'This dataset focuses on live births categorized by five-year age groups of mothers. The data is collected at the "Year" level and covers the territory types "State" and "Region". The dataset includes ...'
$ws.Range("F38").Value = 'OBY03BT02'
$ws.Range("G38").Value = 0.320644736289978
$ws.Range("H38").Value = 0.5692893102279029
$ws.Range("I38").Value = 0.00008220189
$ws.Range("B39").Value = 50
$ws.Range("D39").Value = 40
$ws.Range("E39").Value = 'This dataset focuses on "Housekeeping households" and provides information about the number and types of households in the Czech Republic. The data is collected at the "Census Year" level and covers t...'
$ws.Range("F39").Value = 'SLD039T01'
$ws.Range("G39").Value = 0.3274458646774292
$ws.Range("H39").Value = 0.5686022408218278
$ws.Range("I39").Value = 0.00007843789
$ws.Range("B40").Value = 44
$ws.Range("D40").Value = 37
$ws.Range("E40").Value = 'This dataset focuses on one-family households by type of housekeeping household and regions. The data is collected at the "Census Year" level and covers the territory types "State" and "Region". The d...'
$ws.Range("F40").Value = 'SLD033T05'
$ws.Range("G40").Value = 0.3379700183868408
$ws.Range("H40").Value = 0.5656002811742367
$ws.Range("I40").Value = 0.00016093082
$ws.Range("B41").Value = 46
$ws.Range("D41").Value = 29
$ws.Range("E41").Value = '"This dataset focuses on the main data regarding population and vital statistics, presented in absolute figures. The data is collected at the cumulative quarter level and covers the territory types of...'
$ws.Range("F41").Value = 'OBY01AT01'
$ws.Range("G41").Value = 0.3329526782035828
$ws.Range("H41").Value = 0.5638551449367207
$ws.Range("I41").Value = 0.0009547317
$ws.Range("D42").Value = 51
$ws.Range("E42").Value = '"This dataset focuses on the number of housekeeping household members by regions. The data is collected at the "Census Year" level and covers the territory types "State" and "Region". The dataset incl...'
$ws.Range("F42").Value = 'SLD032T02a'
$ws.Range("G42").Value = 0.3244625926017761
$ws.Range("H42").Value = 0.562759866952105
$ws.Range("I42").Value = 0.000006748052
$ws.Range("B43").ClearContents()
$ws.Range("D43").Value = 57
$ws.Range("E43").Value = 'This dataset focuses on basic information about elections to regional councils by region. The data is collected at the "Year" level and covers the territory types "State" and "Region". The dataset inc...'
$ws.Range("F43").Value = 'VOLKRT2'
$ws.Range("G43").ClearContents()
$ws.Range("H43").Value = 0.5593480275663998
$ws.Range("I43").Value = 0.000004189971
$ws.Range("B44").Value = 55
$ws.Range("D44").Value = 49
$ws.Range("E44").Value = 'This dataset focuses on deaths by sex and basic age groups. The data is collected at the yearly level and covers the territory types of the state and regions. The dataset includes the following column...'
$ws.Range("F44").Value = 'OBY04AT01'
$ws.Range("G44").Value = 0.3213087320327759
$ws.Range("H44").Value = 0.5551534170363642
$ws.Range("I44").Value = 0.000007951138
$ws.Range("B45").Value = 45
$ws.Range("D45").Value = 11
$ws.Range("E45").Value = '"This dataset focuses on population and population change since 1992, expressed as quarterly relative figures. The data is collected at the quarterly level and covers the territory type "State". The d...'
$ws.Range("F45").Value = 'OBY01CQT02'
$ws.Range("G45").Value = 0.3338848352432251
$ws.Range("H45").Value = 0.553739727447636
$ws.Range("I45").Value = 0.007695647
$ws.Range("D46").Value = 41
$ws.Range("E46").Value = 'This dataset focuses on the "Number of women and men by selected age groups in regions." The data is collected at the "Year" level and covers "Region" as the territory type. The dataset includes the f...'
$ws.Range("F46").Value = 'WGEN01G1VN'
$ws.Range("H46").Value = 0.5536686474604338
$ws.Range("I46").Value = 0.000041986306
$ws.Range("B47").ClearContents()
$ws.Range("D47").Value = 55
$ws.Range("E47").Value = 'This dataset focuses on basic information regarding the election of the President of the Czech Republic by region. The data is collected at the "Year" level and covers the territory types "State" and ...'
$ws.Range("F47").Value = 'VOLPRT2'
$ws.Range("G47").ClearContents()
$ws.Range("H47").Value = 0.5526561475800364
$ws.Range("I47").Value = 0.000005862814
$ws.Range("B48").ClearContents()
$ws.Range("D48").Value = 42
$ws.Range("E48").Value = 'This dataset focuses on occupied houses categorized by the period of construction or reconstruction and by regions. The data is collected at the "Census Year" level and covers the territory types "Sta...'
$ws.Range("F48").Value = 'SLD020T02'
$ws.Range("G48").ClearContents()
$ws.Range("H48").Value = 0.5513126466934043
$ws.Range("I48").Value = 0.00003883111
$ws.Range("B49").Value = 59
$ws.Range("D49").Value = 38
$ws.Range("E49").Value = '"This dataset focuses on the total population as at 1. 1., 1. 7., and 31. 12., categorized by sex. The data is collected at the yearly level and covers the territory types of the state and regions. Th...'
$ws.Range("F49").Value = 'OBY02AT01'
$ws.Range("G49").Value = 0.3194942474365234
$ws.Range("H49").Value = 0.551102191997754
$ws.Range("I49").Value = 0.00012339458
$ws.Range("B50").Value = 58
$ws.Range("D50").Value = 59
$ws.Range("E50").Value = 'This dataset focuses on the population by housing arrangements and regions. The data is collected at the "Census Year" level and covers the territory types "State" and "Region". The dataset includes t...'
$ws.Range("F50").Value = 'SLD013T02'
$ws.Range("G50").Value = 0.3204777240753174
$ws.Range("H50").Value = 0.550091251866292
$ws.Range("I50").Value = 0.000003314543
$ws.Range("D51").Value = 52
$ws.Range("D53").Value = 56
$ws.Range("D54").Value = 54
$ws.Range("D55").Value = 48
$ws.Range("D56").Value = 45
$ws.Range("D58").Value = 46
$ws.Range("D59").Value = 58
$ws.Range("D60").Value = 50
$ws.Range("D61").Value = 53
$ws.Range("B62").Value = 52
$ws.Range("E62").Value = 'This dataset focuses on "Population change - relative figures per 1,000 population". The data is collected at the "Cumulative Quarter" level and covers the territory types "State" and "Region". The da...'
$ws.Range("F62").Value = 'OBY01PDT03'
$ws.Range("G62").Value = 0.3243331909179688
$ws.Range("B63").Value = 60
$ws.Range("E63").Value = '"This dataset focuses on occupied dwellings categorized by the tenure status of the dwelling and the region. The data is collected at the ''Census Year'' level and covers the territory types ''State'' and...'
$ws.Range("F63").Value = 'SLD024T02'
$ws.Range("G63").Value = 0.3193652629852295
$ws.Range("B64").Value = 57
$ws.Range("E64").Value = '"This dataset focuses on population change in absolute figures. The data is collected at the cumulative quarter level and covers state and regional territory types. The dataset includes the following ...'
$ws.Range("F64").Value = 'OBY01PDT02'
$ws.Range("G64").Value = 0.320556640625
$ws.Range("B65").Value = 54
$ws.Range("E65").Value = '"This dataset focuses on the main data regarding population and vital statistics, specifically relative figures. The data is collected at the ''Cumulative Quarter'' level and covers the territory types ...'
$ws.Range("F65").Value = 'OBY01AT02'
$ws.Range("G65").Value = 0.3217225074768066
$ws.Range("B66").Value = 53
$ws.Range("E66").Value = '"This dataset focuses on one-family households categorized by the number of dependent children and regions. The data is collected at the "Census Year" level and covers "State" and "Region" territory t...'
$ws.Range("F66").Value = 'SLD033T04a'
$ws.Range("G66").Value = 0.3242666125297546
